# Daily attendance processing - 2026-01-29 16:58:32
# Swap the order of the two comma-separated "Recorded By" entries in
# column G whenever "dnasr281@gmail.com" appears as the second (last)
# entry alongside exactly one other recorder, e.g.
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "admin@admin.com, dnasr281@gmail.com"   -> "dnasr281@gmail.com, admin@admin.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*, dnasr281@gmail.com") {
        $parts = $val -split ", "
        if ($parts.Count -eq 2) {
            $newVal = $parts[1] + ", " + $parts[0]
            $cell.Value2 = $newVal
        }
    }
}
